$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1 — copy formatting from the neighboring
# header cell (G1) so it gets the same bold/border/centered style, then
# set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Save column data for the two existing rows.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
